$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: 'Datos actualizados a 31 de Mayo de 2020 a las 12:05' -> 'Datos actualizados a 31 de Mayo de 2020 a las 13:10'
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 31 de Mayo de 2020 a las 13:10'

# Row 12: 'India' -> 'India'
$ws.Cells.Item(12, 2).Value = 182990
$ws.Cells.Item(12, 3).Value = 1163
$ws.Cells.Item(12, 4).Value = 87099
$ws.Cells.Item(12, 5).Value = 90703
$ws.Cells.Item(12, 7).Value = 3
$ws.Cells.Item(12, 8).Value = 5188

# Row 15: 'Iran' -> 'Iran'
$ws.Cells.Item(15, 2).Value = 151466
$ws.Cells.Item(15, 3).Value = 2516
$ws.Cells.Item(15, 4).Value = 118848
$ws.Cells.Item(15, 5).Value = 24821
$ws.Cells.Item(15, 7).Value = 63
$ws.Cells.Item(15, 8).Value = 7797

# Row 23: 'Catar' -> 'Catar'
$ws.Cells.Item(23, 2).Value = 56910
$ws.Cells.Item(23, 3).Value = 1648
$ws.Cells.Item(23, 4).Value = 30290
$ws.Cells.Item(23, 5).Value = 26582
$ws.Cells.Item(23, 7).Value = 2
$ws.Cells.Item(23, 8).Value = 38

# Row 33: 'Suiza' -> 'Suiza'
$ws.Cells.Item(33, 2).Value = 30862
$ws.Cells.Item(33, 3).Value = 17
$ws.Cells.Item(33, 5).Value = 542
$ws.Cells.Item(33, 7).Value = 1
$ws.Cells.Item(33, 8).Value = 1920

# Row 35: 'Indonesia' -> 'Kuwait'
$ws.Cells.Item(35, 1).Value = 'Kuwait'
$ws.Cells.Item(35, 2).Value = 27043
$ws.Cells.Item(35, 3).Value = 851
$ws.Cells.Item(35, 4).Value = 11386
$ws.Cells.Item(35, 5).Value = 15445
$ws.Cells.Item(35, 7).Value = 7
$ws.Cells.Item(35, 8).Value = 212

# Row 36: 'Kuwait' -> 'Indonesia'
$ws.Cells.Item(36, 1).Value = 'Indonesia'
$ws.Cells.Item(36, 2).Value = 26473
$ws.Cells.Item(36, 3).Value = 700
$ws.Cells.Item(36, 4).Value = 7308
$ws.Cells.Item(36, 5).Value = 17552
$ws.Cells.Item(36, 7).Value = 40
$ws.Cells.Item(36, 8).Value = 1613

# Row 46: 'Austria' -> 'Austria'
$ws.Cells.Item(46, 2).Value = 16731
$ws.Cells.Item(46, 3).Value = 46
$ws.Cells.Item(46, 4).Value = 15593
$ws.Cells.Item(46, 5).Value = 470

# Row 52: 'Serbia' -> 'Oman'
$ws.Cells.Item(52, 1).Value = 'Oman'
$ws.Cells.Item(52, 2).Value = 11437
$ws.Cells.Item(52, 3).Value = 1014
$ws.Cells.Item(52, 4).Value = 2396
$ws.Cells.Item(52, 5).Value = 8997
$ws.Cells.Item(52, 7).Value = 2
$ws.Cells.Item(52, 8).Value = 44

# Row 53: 'Barein' -> 'Serbia'
$ws.Cells.Item(53, 1).Value = 'Serbia'
$ws.Cells.Item(53, 2).Value = 11381
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 6606
$ws.Cells.Item(53, 5).Value = 4533
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 242

# Row 54: 'Kazajistan' -> 'Barein'
$ws.Cells.Item(54, 1).Value = 'Barein'
$ws.Cells.Item(54, 2).Value = 11288
$ws.Cells.Item(54, 3).Value = 495
$ws.Cells.Item(54, 4).Value = 6673
$ws.Cells.Item(54, 5).Value = 4597
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 18

# Row 55: 'Oman' -> 'Kazajistan'
$ws.Cells.Item(55, 1).Value = 'Kazajistan'
$ws.Cells.Item(55, 2).Value = 10858
$ws.Cells.Item(55, 3).Value = 476
$ws.Cells.Item(55, 4).Value = 5220
$ws.Cells.Item(55, 5).Value = 5600
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 38

# Row 78: 'Uzbekistan' -> 'Senegal'
$ws.Cells.Item(78, 1).Value = 'Senegal'
$ws.Cells.Item(78, 2).Value = 3645
$ws.Cells.Item(78, 3).Value = 110
$ws.Cells.Item(78, 4).Value = 1801
$ws.Cells.Item(78, 5).Value = 1802
$ws.Cells.Item(78, 8).Value = 42

# Row 79: 'Senegal' -> 'Uzbekistan'
$ws.Cells.Item(79, 1).Value = 'Uzbekistan'
$ws.Cells.Item(79, 2).Value = 3554
$ws.Cells.Item(79, 3).Value = 8
$ws.Cells.Item(79, 4).Value = 2783
$ws.Cells.Item(79, 5).Value = 757
$ws.Cells.Item(79, 8).Value = 14

# Row 88: 'Bosnia y Herzegovina' -> 'Bosnia y Herzegovina'
$ws.Cells.Item(88, 2).Value = 2510
$ws.Cells.Item(88, 3).Value = 16
$ws.Cells.Item(88, 4).Value = 1862
$ws.Cells.Item(88, 5).Value = 495

# Row 102: 'Eslovaquia' -> 'Nepal'
$ws.Cells.Item(102, 1).Value = 'Nepal'
$ws.Cells.Item(102, 2).Value = 1567
$ws.Cells.Item(102, 3).Value = 166
$ws.Cells.Item(102, 4).Value = 220
$ws.Cells.Item(102, 5).Value = 1339
$ws.Cells.Item(102, 7).Value = 2
$ws.Cells.Item(102, 8).Value = 8

# Row 103: 'Nueva Zelanda' -> 'Eslovaquia'
$ws.Cells.Item(103, 1).Value = 'Eslovaquia'
$ws.Cells.Item(103, 2).Value = 1521
$ws.Cells.Item(103, 4).Value = 1366
$ws.Cells.Item(103, 5).Value = 127
$ws.Cells.Item(103, 8).Value = 28

# Row 104: 'Eslovenia' -> 'Nueva Zelanda'
$ws.Cells.Item(104, 1).Value = 'Nueva Zelanda'
$ws.Cells.Item(104, 2).Value = 1504
$ws.Cells.Item(104, 4).Value = 1481
$ws.Cells.Item(104, 5).Value = 1
$ws.Cells.Item(104, 8).Value = 22

# Row 105: 'Venezuela' -> 'Eslovenia'
$ws.Cells.Item(105, 1).Value = 'Eslovenia'
$ws.Cells.Item(105, 2).Value = 1473
$ws.Cells.Item(105, 4).Value = 1358
$ws.Cells.Item(105, 5).Value = 7
$ws.Cells.Item(105, 8).Value = 108

# Row 106: 'Nepal' -> 'Venezuela'
$ws.Cells.Item(106, 1).Value = 'Venezuela'
$ws.Cells.Item(106, 2).Value = 1459
$ws.Cells.Item(106, 4).Value = 302
$ws.Cells.Item(106, 5).Value = 1143
$ws.Cells.Item(106, 8).Value = 14

# Row 134: 'Malta' -> 'Malta'
$ws.Cells.Item(134, 4).Value = 534
$ws.Cells.Item(134, 5).Value = 75

# Row 200: 'Belice' -> 'Santa Lucia'
$ws.Cells.Item(200, 1).Value = 'Santa Lucia'
$ws.Cells.Item(200, 4).Value = 18
$ws.Cells.Item(200, 8).Value = 0

# Row 201: 'Santa Lucia' -> 'Belice'
$ws.Cells.Item(201, 1).Value = 'Belice'
$ws.Cells.Item(201, 4).Value = 16
$ws.Cells.Item(201, 8).Value = 2

# Row 213: 'Islas Virgenes Britanicas' -> 'Papua Nueva Guinea'
$ws.Cells.Item(213, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 8).Value = 0

# Row 214: 'Papua Nueva Guinea' -> 'Islas Virgenes Britanicas'
$ws.Cells.Item(214, 1).Value = 'Islas Virgenes Britanicas'
$ws.Cells.Item(214, 4).Value = 7

# Row 218: 'Lesoto' -> 'Lesoto'
$ws.Cells.Item(218, 4).Value = 0
